$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.62"
$ws.Range("E2").Value = "'-0.23%"
$ws.Range("G2").Value = "'15"
$ws.Range("D3").Value = "'29.36"
$ws.Range("E3").Value = "'-0.48%"
$ws.Range("G3").Value = "'15"
$ws.Range("D4").Value = "'5.304"
$ws.Range("E4").Value = "'2.35%"
$ws.Range("G4").Value = "'15"
$ws.Range("D5").Value = "'0.05728"
$ws.Range("E5").Value = "'0.41%"
$ws.Range("G5").Value = "'15"
$ws.Range("D6").Value = "'6.636"
$ws.Range("E6").Value = "'0.72%"
$ws.Range("G6").Value = "'15"
$ws.Range("D7").Value = "'3.193"
$ws.Range("E7").Value = "'4.52%"
$ws.Range("G7").Value = "'15"
$ws.Range("D8").Value = "'0.8570"
$ws.Range("E8").Value = "'-0.20%"
$ws.Range("G8").Value = "'15"
$ws.Range("D9").Value = "'0.8545"
$ws.Range("E9").Value = "'-2.88%"
$ws.Range("G9").Value = "'15"
$ws.Range("D10").Value = "'0.1378"
$ws.Range("E10").Value = "'0.71%"
$ws.Range("G10").Value = "'15"
$ws.Range("D11").Value = "'0.07093"
$ws.Range("E11").Value = "'-0.12%"
$ws.Range("G11").Value = "'15"
$ws.Range("D12").Value = "'0.03184"
$ws.Range("E12").Value = "'11.10%"
$ws.Range("G12").Value = "'15"
$ws.Range("D13").Value = "'0.09345"
$ws.Range("E13").Value = "'-0.47%"
$ws.Range("G13").Value = "'15"
$ws.Range("D14").Value = "'0.001552"
$ws.Range("E14").Value = "'2.04%"
$ws.Range("G14").Value = "'15"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0005961"
$ws.Range("E15").Value = "'-94.22%"
$ws.Range("G15").Value = "'15"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006050"
$ws.Range("E16").Value = "'-0.75%"
$ws.Range("G16").Value = "'15"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.516"
$ws.Range("E17").Value = "'0.94%"
$ws.Range("G17").Value = "'15"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.194"
$ws.Range("E18").Value = "'0.39%"
$ws.Range("G18").Value = "'15"
$ws.Range("D19").Value = "'0.3163"
$ws.Range("E19").Value = "'0.63%"
$ws.Range("G19").Value = "'15"
$ws.Range("D20").Value = "'0.03341"
$ws.Range("E20").Value = "'1.57%"
$ws.Range("G20").Value = "'15"
$ws.Range("D21").Value = "'0.1305"
$ws.Range("E21").Value = "'0.33%"
$ws.Range("G21").Value = "'15"
$ws.Range("D22").Value = "'3.486"
$ws.Range("E22").Value = "'0.49%"
$ws.Range("G22").Value = "'15"
$ws.Range("D23").Value = "'0.04135"
$ws.Range("E23").Value = "'-0.50%"
$ws.Range("G23").Value = "'15"
$ws.Range("D24").Value = "'0.1411"
$ws.Range("E24").Value = "'2.22%"
$ws.Range("G24").Value = "'15"
$ws.Range("D25").Value = "'0.001225"
$ws.Range("E25").Value = "'0.20%"
$ws.Range("G25").Value = "'15"
$ws.Range("D26").Value = "'0.004157"
$ws.Range("E26").Value = "'-18.31%"
$ws.Range("G26").Value = "'15"
$ws.Range("D27").Value = "'0.0001201"
$ws.Range("E27").Value = "'-0.80%"
$ws.Range("G27").Value = "'15"
$ws.Range("D28").Value = "'0.0001452"
$ws.Range("E28").Value = "'-25.11%"
$ws.Range("G28").Value = "'15"
$ws.Range("G29").Value = "'15"
$ws.Range("G30").Value = "'15"
$ws.Range("G31").Value = "'15"
$ws.Range("G32").Value = "'15"
$ws.Range("G33").Value = "'15"
$ws.Range("G34").Value = "'15"
$ws.Range("G35").Value = "'15"
$ws.Range("G36").Value = "'15"
$ws.Range("G37").Value = "'15"
$ws.Range("G38").Value = "'15"
$ws.Range("G39").Value = "'15"
$ws.Range("D40").Value = "'0.03757"
$ws.Range("E40").Value = "'0.25%"
$ws.Range("G40").Value = "'15"
$ws.Range("D41").Value = "'0.1067"
$ws.Range("E41").Value = "'-0.45%"
$ws.Range("G41").Value = "'15"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.003686"
$ws.Range("E42").Value = "'-36.01%"
$ws.Range("G42").Value = "'15"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002419"
$ws.Range("E43").Value = "'-4.77%"
$ws.Range("G43").Value = "'15"
$ws.Range("D44").Value = "'0.009943"
$ws.Range("E44").Value = "'5.69%"
$ws.Range("G44").Value = "'15"
$ws.Range("D45").Value = "'0.00005289"
$ws.Range("E45").Value = "'3.53%"
$ws.Range("G45").Value = "'15"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'0.19%"
$ws.Range("G46").Value = "'15"
$ws.Range("D47").Value = "'0.08104"
$ws.Range("E47").Value = "'14.15%"
$ws.Range("G47").Value = "'15"
$ws.Range("D48").Value = "'0.002195"
$ws.Range("E48").Value = "'-17.58%"
$ws.Range("G48").Value = "'15"
$ws.Range("D49").Value = "'0.00002104"
$ws.Range("E49").Value = "'0.19%"
$ws.Range("G49").Value = "'15"
$ws.Range("D50").Value = "'0.0002004"
$ws.Range("E50").Value = "'0.19%"
$ws.Range("G50").Value = "'15"
$ws.Range("G51").Value = "'15"
